# Apply updated NATA airtoxics data to the "Means" and "Standard Deviations" sheets.

$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# --- Means sheet ---
# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 24
$wsMeans.Range("D9").Value = 20
$wsMeans.Range("E9").Value = 20
$wsMeans.Range("F9").Value = 20
$wsMeans.Range("G9").Value = 20

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.29
$wsMeans.Range("D10").Value = 0.3
$wsMeans.Range("E10").Value = 0.27
$wsMeans.Range("F10").Value = 0.26
$wsMeans.Range("G10").Value = 0.27

# --- Standard Deviations sheet ---
# Row 9: Total Cancer Risk (per million) SD
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 5.6
$wsSD.Range("D9").Value = 0
$wsSD.Range("E9").Value = 0
$wsSD.Range("F9").Value = 0
$wsSD.Range("G9").Value = 0

# Row 10: Total Respiratory (hazard quotient) SD
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.086
$wsSD.Range("D10").Value = 0
$wsSD.Range("E10").Value = 0.046
$wsSD.Range("F10").Value = 0.046
$wsSD.Range("G10").Value = 0.043

$wb.Save()
